$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B-G (ME, MAE, MSE, RMSE, SE, N)
$data = @{
    2  = @(-0.100195032778358, 1.879008151558315, 15.52109306684352, 3.939681848429327, 4.026922141332304, 23)
    3  = @(0.1448133846575885, 1.786984302705635, 10.36886832397854, 3.220072720293524, 3.292514909155644, 22)
    4  = @(-0.5106375609785742, 1.177604939791601, 4.283545554090664, 2.069672813294571, 2.055221296400692, 21)
    5  = @(0.07283602206394024, 0.8438801149343232, 1.589728322734128, 1.260844289646476, 1.291438706748583, 20)
    6  = @(0.009315935962682875, 0.8519364303311395, 1.743403430867843, 1.320380032743544, 1.356527761928607, 19)
    7  = @(0.01805680939183702, 0.6865640439233759, 0.7225809383130525, 0.8500476094390552, 0.8744944091926862, 18)
    8  = @(0.06286173468651945, 0.8035601421598704, 1.081400171720104, 1.039903924273827, 1.069948180945846, 17)
    9  = @(0.2080312578962194, 0.6959573112162807, 0.6635098639986297, 0.8145611481028479, 0.813376737467042, 16)
    10 = @(0.1458528612741254, 0.6988170516431337, 0.9094933485271182, 0.9536736069154468, 0.9755329814582326, 15)
    11 = @(0.2106889511877471, 0.5147298635682541, 0.3550355181284862, 0.5958485697964595, 0.578395630879239, 14)
}

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    $ws.Cells.Item($rowNum, 2).Value = $values[0]
    $ws.Cells.Item($rowNum, 3).Value = $values[1]
    $ws.Cells.Item($rowNum, 4).Value = $values[2]
    $ws.Cells.Item($rowNum, 5).Value = $values[3]
    $ws.Cells.Item($rowNum, 6).Value = $values[4]
    $ws.Cells.Item($rowNum, 7).Value = $values[5]
}
